$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44315
$ws.Cells.Item(2, 11).Value = 4000
$ws.Cells.Item(2, 12).Value = 4000
$ws.Cells.Item(2, 13).Value = 4000
$ws.Cells.Item(2, 16).Value = 4000

$ws.Cells.Item(3, 4).Value = 44956
$ws.Cells.Item(3, 10).Value = 40

$ws.Cells.Item(4, 4).Value = 44508
$ws.Cells.Item(4, 10).Value = 30

$ws.Cells.Item(5, 4).Value = 44313
$ws.Cells.Item(5, 10).Value = 20
$ws.Cells.Item(5, 11).Value = 4000
$ws.Cells.Item(5, 12).Value = 4000
$ws.Cells.Item(5, 13).Value = 4000
$ws.Cells.Item(5, 16).Value = 4000

$ws.Cells.Item(6, 4).Value = 44498
$ws.Cells.Item(6, 10).Value = 40

$ws.Cells.Item(7, 4).Value = 44749
$ws.Cells.Item(7, 10).Value = 65
$ws.Cells.Item(7, 11).Value = 6000
$ws.Cells.Item(7, 12).Value = 6000
$ws.Cells.Item(7, 13).Value = 6000
$ws.Cells.Item(7, 16).Value = 6000

$ws.Cells.Item(8, 4).Value = 44679
$ws.Cells.Item(8, 10).Value = 50
$ws.Cells.Item(8, 11).Value = 5000
$ws.Cells.Item(8, 12).Value = 5000
$ws.Cells.Item(8, 13).Value = 5000
$ws.Cells.Item(8, 16).Value = 5000

$ws.Cells.Item(9, 4).Value = 44509

$ws.Cells.Item(10, 4).Value = 44504
$ws.Cells.Item(10, 10).Value = 55
$ws.Cells.Item(10, 11).Value = 4000
$ws.Cells.Item(10, 12).Value = 4000
$ws.Cells.Item(10, 13).Value = 4000
$ws.Cells.Item(10, 16).Value = 4000

$ws.Cells.Item(11, 4).Value = 44301
$ws.Cells.Item(11, 10).Value = 40
$ws.Cells.Item(11, 11).Value = 3000
$ws.Cells.Item(11, 12).Value = 3000
$ws.Cells.Item(11, 13).Value = 3000
$ws.Cells.Item(11, 16).Value = 3000

$ws.Cells.Item(12, 4).Value = 44649
$ws.Cells.Item(12, 10).Value = 20

$ws.Cells.Item(13, 4).Value = 44777
$ws.Cells.Item(13, 10).Value = 25
$ws.Cells.Item(13, 11).Value = 5000
$ws.Cells.Item(13, 12).Value = 5000
$ws.Cells.Item(13, 13).Value = 5000
$ws.Cells.Item(13, 16).Value = 5000

$ws.Cells.Item(14, 4).Value = 44966
$ws.Cells.Item(14, 10).Value = 40

$ws.Cells.Item(15, 4).Value = 44365
$ws.Cells.Item(15, 10).Value = 55
$ws.Cells.Item(15, 11).Value = 5000
$ws.Cells.Item(15, 12).Value = 5000
$ws.Cells.Item(15, 13).Value = 5000
$ws.Cells.Item(15, 16).Value = 5000

$ws.Cells.Item(16, 4).Value = 44312
$ws.Cells.Item(16, 10).Value = 50
$ws.Cells.Item(16, 11).Value = 4000
$ws.Cells.Item(16, 12).Value = 4000
$ws.Cells.Item(16, 13).Value = 4000
$ws.Cells.Item(16, 16).Value = 4000

$ws.Cells.Item(17, 4).Value = 44390
$ws.Cells.Item(17, 11).Value = 6000
$ws.Cells.Item(17, 12).Value = 6000
$ws.Cells.Item(17, 13).Value = 6000
$ws.Cells.Item(17, 16).Value = 6000

$ws.Cells.Item(18, 4).Value = 44259
$ws.Cells.Item(18, 10).Value = 30

$ws.Cells.Item(19, 4).Value = 44316

$ws.Cells.Item(20, 4).Value = 44176
$ws.Cells.Item(20, 10).Value = 10
$ws.Cells.Item(20, 11).Value = 4000
$ws.Cells.Item(20, 12).Value = 4000
$ws.Cells.Item(20, 13).Value = 4000
$ws.Cells.Item(20, 16).Value = 4000

$ws.Cells.Item(21, 4).Value = 44781
$ws.Cells.Item(21, 11).Value = 5000
$ws.Cells.Item(21, 12).Value = 5000
$ws.Cells.Item(21, 13).Value = 5000
$ws.Cells.Item(21, 16).Value = 5000

$ws.Cells.Item(22, 4).Value = 44291
$ws.Cells.Item(22, 10).Value = 35
$ws.Cells.Item(22, 11).Value = 4000
$ws.Cells.Item(22, 12).Value = 4000
$ws.Cells.Item(22, 13).Value = 4000
$ws.Cells.Item(22, 16).Value = 4000

$ws.Cells.Item(23, 4).Value = 44280
$ws.Cells.Item(23, 10).Value = 55
$ws.Cells.Item(23, 11).Value = 4000
$ws.Cells.Item(23, 12).Value = 4000
$ws.Cells.Item(23, 13).Value = 4000
$ws.Cells.Item(23, 16).Value = 4000

$ws.Cells.Item(24, 4).Value = 44680
$ws.Cells.Item(24, 10).Value = 20

$ws.Cells.Item(25, 4).Value = 44959
$ws.Cells.Item(25, 11).Value = 5000
$ws.Cells.Item(25, 12).Value = 5000
$ws.Cells.Item(25, 13).Value = 5000
$ws.Cells.Item(25, 16).Value = 5000

$ws.Cells.Item(26, 4).Value = 44497
$ws.Cells.Item(26, 10).Value = 20
$ws.Cells.Item(26, 11).Value = 4000
$ws.Cells.Item(26, 12).Value = 4000
$ws.Cells.Item(26, 13).Value = 4000
$ws.Cells.Item(26, 16).Value = 4000

$ws.Cells.Item(27, 4).Value = 44656
$ws.Cells.Item(27, 10).Value = 85
$ws.Cells.Item(27, 11).Value = 5000
$ws.Cells.Item(27, 12).Value = 5000
$ws.Cells.Item(27, 13).Value = 5000
$ws.Cells.Item(27, 16).Value = 5000
